$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Identity")

# Step 1: insert 39 new rows above the old bottom border row (180 -> 219),
# copying formatting from row 179 (the last pre-existing placeholder row)
# so the newly inserted rows pick up the same borders/styles.
$ws.Rows("180:218").Insert()
$fmtSrc = $ws.Range("A179:J179")
$fmtDst = $ws.Range("A180:J218")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 2: fill in the new test-plan rows (column A matches the order the
# author typed the test titles; column B follows the authors ticket-number
# entry order so the shared-string table builds up the same way).

# -- column A, rows 172-182 --
$ws.Range("A172").Value2 = 'User cannot be created if username is not supplied'
$ws.Range("A173").Value2 = 'User cannot be created if username is not supplied'
$ws.Range("A174").Value2 = 'User cannot be created if email is not supplied'
$ws.Range("A175").Value2 = 'User cannot be created if email is invalid'
$ws.Range("A176").Value2 = 'User cannot be created if display name is not supplied'
$ws.Range("A177").Value2 = 'User cannot be created if custom metadata validations are not met'
$ws.Range("A178").Value2 = 'User cannot be created if password is not supplied'
$ws.Range("A179").Value2 = 'User cannot be created if passwords do not metch'
$ws.Range("A180").Value2 = 'User cannot be created if no role is supplied'
$ws.Range("A181").Value2 = 'User cannot be created if requesting user doesn''t have the proper access rights'
$ws.Range("A182").Value2 = 'User can be created'

# -- column B, rows 182, 181, then 172-180 --
$ws.Range("B182").Value2 = 'IT-237'
$ws.Range("B181").Value2 = 'IT-236'
$ws.Range("B172").Value2 = 'IT-227'
$ws.Range("B173").Value2 = 'IT-228'
$ws.Range("B174").Value2 = 'IT-229'
$ws.Range("B175").Value2 = 'IT-230'
$ws.Range("B176").Value2 = 'IT-231'
$ws.Range("B177").Value2 = 'IT-232'
$ws.Range("B178").Value2 = 'IT-233'
$ws.Range("B179").Value2 = 'IT-234'
$ws.Range("B180").Value2 = 'IT-235'

# -- columns D-G, rows 172-182 --
$ws.Range("D172").Value2 = 'Identity'; $ws.Range("E172").Value2 = 'User'; $ws.Range("F172").Value2 = 'General'; $ws.Range("G172").Value2 = 'Validation'
$ws.Range("D173").Value2 = 'Identity'; $ws.Range("E173").Value2 = 'User'; $ws.Range("F173").Value2 = 'General'; $ws.Range("G173").Value2 = 'Validation'
$ws.Range("D174").Value2 = 'Identity'; $ws.Range("E174").Value2 = 'User'; $ws.Range("F174").Value2 = 'General'; $ws.Range("G174").Value2 = 'Validation'
$ws.Range("D175").Value2 = 'Identity'; $ws.Range("E175").Value2 = 'User'; $ws.Range("F175").Value2 = 'General'; $ws.Range("G175").Value2 = 'Validation'
$ws.Range("D176").Value2 = 'Identity'; $ws.Range("E176").Value2 = 'User'; $ws.Range("F176").Value2 = 'General'; $ws.Range("G176").Value2 = 'Validation'
$ws.Range("D177").Value2 = 'Identity'; $ws.Range("E177").Value2 = 'User'; $ws.Range("F177").Value2 = 'General'; $ws.Range("G177").Value2 = 'Validation'
$ws.Range("D178").Value2 = 'Identity'; $ws.Range("E178").Value2 = 'User'; $ws.Range("F178").Value2 = 'General'; $ws.Range("G178").Value2 = 'Validation'
$ws.Range("D179").Value2 = 'Identity'; $ws.Range("E179").Value2 = 'User'; $ws.Range("F179").Value2 = 'General'; $ws.Range("G179").Value2 = 'Validation'
$ws.Range("D180").Value2 = 'Identity'; $ws.Range("E180").Value2 = 'User'; $ws.Range("F180").Value2 = 'General'; $ws.Range("G180").Value2 = 'Validation'
$ws.Range("D181").Value2 = 'Identity'; $ws.Range("E181").Value2 = 'User'; $ws.Range("F181").Value2 = 'General'; $ws.Range("G181").Value2 = 'Validation'
$ws.Range("D182").Value2 = 'Identity'; $ws.Range("E182").Value2 = 'User'; $ws.Range("F182").Value2 = 'General'; $ws.Range("G182").Value2 = 'Business Logic'

# -- column B, rows 183-194 --
$ws.Range("B183").Value2 = 'IT-250'
$ws.Range("B184").Value2 = 'IT-251'
$ws.Range("B185").Value2 = 'IT-252'
$ws.Range("B186").Value2 = 'IT-253'
$ws.Range("B187").Value2 = 'IT-254'
$ws.Range("B188").Value2 = 'IT-255'
$ws.Range("B189").Value2 = 'IT-256'
$ws.Range("B190").Value2 = 'IT-257'
$ws.Range("B191").Value2 = 'IT-258'
$ws.Range("B192").Value2 = 'IT-259'
$ws.Range("B193").Value2 = 'IT-260'
$ws.Range("B194").Value2 = 'IT-261'

# -- column A, rows 183-194 --
$ws.Range("A183").Value2 = 'User cannot register if username is not provided'
$ws.Range("A184").Value2 = 'User cannot register if username is already taken'
$ws.Range("A185").Value2 = 'User cannot register if tenant does not exist'
$ws.Range("A186").Value2 = 'User cannot register if tenant is not provided'
$ws.Range("A187").Value2 = 'User cannot register if email is not provided'
$ws.Range("A188").Value2 = 'User cannot register if email is not valid'
$ws.Range("A189").Value2 = 'User cannot register if tenant is display name is not provided'
$ws.Range("A190").Value2 = 'User cannot register if metadata custom validators are not met'
$ws.Range("A191").Value2 = 'User cannot register if password is not provided'
$ws.Range("A192").Value2 = 'User cannot register if passwords don''t match'
$ws.Range("A193").Value2 = 'User cannot register if passwords custom policies are not met'
$ws.Range("A194").Value2 = 'User can register'

# -- columns D-G, rows 183-194 --
$ws.Range("D183").Value2 = 'Identity'; $ws.Range("E183").Value2 = 'User'; $ws.Range("F183").Value2 = 'General'; $ws.Range("G183").Value2 = 'Validation'
$ws.Range("D184").Value2 = 'Identity'; $ws.Range("E184").Value2 = 'User'; $ws.Range("F184").Value2 = 'General'; $ws.Range("G184").Value2 = 'Validation'
$ws.Range("D185").Value2 = 'Identity'; $ws.Range("E185").Value2 = 'User'; $ws.Range("F185").Value2 = 'General'; $ws.Range("G185").Value2 = 'Validation'
$ws.Range("D186").Value2 = 'Identity'; $ws.Range("E186").Value2 = 'User'; $ws.Range("F186").Value2 = 'General'; $ws.Range("G186").Value2 = 'Validation'
$ws.Range("D187").Value2 = 'Identity'; $ws.Range("E187").Value2 = 'User'; $ws.Range("F187").Value2 = 'General'; $ws.Range("G187").Value2 = 'Validation'
$ws.Range("D188").Value2 = 'Identity'; $ws.Range("E188").Value2 = 'User'; $ws.Range("F188").Value2 = 'General'; $ws.Range("G188").Value2 = 'Validation'
$ws.Range("D189").Value2 = 'Identity'; $ws.Range("E189").Value2 = 'User'; $ws.Range("F189").Value2 = 'General'; $ws.Range("G189").Value2 = 'Validation'
$ws.Range("D190").Value2 = 'Identity'; $ws.Range("E190").Value2 = 'User'; $ws.Range("F190").Value2 = 'General'; $ws.Range("G190").Value2 = 'Validation'
$ws.Range("D191").Value2 = 'Identity'; $ws.Range("E191").Value2 = 'User'; $ws.Range("F191").Value2 = 'General'; $ws.Range("G191").Value2 = 'Validation'
$ws.Range("D192").Value2 = 'Identity'; $ws.Range("E192").Value2 = 'User'; $ws.Range("F192").Value2 = 'General'; $ws.Range("G192").Value2 = 'Validation'
$ws.Range("D193").Value2 = 'Identity'; $ws.Range("E193").Value2 = 'User'; $ws.Range("F193").Value2 = 'General'; $ws.Range("G193").Value2 = 'Validation'
$ws.Range("D194").Value2 = 'Identity'; $ws.Range("E194").Value2 = 'User'; $ws.Range("F194").Value2 = 'General'; $ws.Range("G194").Value2 = 'Business Logic'

# Step 3: update dimension-driving view state (frozen pane anchor + active cell)
$ws.Activate() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 166
$ws.Range("A173").Select() | Out-Null
